# Update the "10 min frame" forecast report data (rows 2-20, columns B-J)
# with refreshed figures, per the source report regeneration.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = 732
$ws.Cells.Item(2, 4).Value = 739
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 3
$ws.Cells.Item(2, 7).Value = 13
$ws.Cells.Item(2, 9).Value = 806
$ws.Cells.Item(2, 10).Value = -8.31265508684863

# Row 3
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 268
$ws.Cells.Item(3, 4).Value = 269
$ws.Cells.Item(3, 7).Value = 4
$ws.Cells.Item(3, 9).Value = 338
$ws.Cells.Item(3, 10).Value = -20.41420118343196

# Row 4
$ws.Cells.Item(4, 3).Value = 11
$ws.Cells.Item(4, 4).Value = 11
$ws.Cells.Item(4, 9).Value = 7
$ws.Cells.Item(4, 10).Value = 57.14285714285714

# Row 5
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 240
$ws.Cells.Item(5, 4).Value = 256
$ws.Cells.Item(5, 5).Value = 13
$ws.Cells.Item(5, 6).Value = 3
$ws.Cells.Item(5, 9).Value = 144
$ws.Cells.Item(5, 10).Value = 77.77777777777777

# Row 6
$ws.Cells.Item(6, 3).Value = 56
$ws.Cells.Item(6, 4).Value = 57
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 9).Value = 48
$ws.Cells.Item(6, 10).Value = 18.75

# Row 7
$ws.Cells.Item(7, 3).Value = 107
$ws.Cells.Item(7, 4).Value = 108
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 4
$ws.Cells.Item(7, 9).Value = 110
$ws.Cells.Item(7, 10).Value = -1.818181818181819

# Row 8
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 303
$ws.Cells.Item(8, 4).Value = 311
$ws.Cells.Item(8, 5).Value = 7
$ws.Cells.Item(8, 6).Value = 2
$ws.Cells.Item(8, 7).Value = 9
$ws.Cells.Item(8, 10).Value = 972.4137931034483

# Row 9
$ws.Cells.Item(9, 3).Value = 39
$ws.Cells.Item(9, 4).Value = 42
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 10).Value = -40.84507042253522

# Row 10
$ws.Cells.Item(10, 2).Value = 7
$ws.Cells.Item(10, 3).Value = 471
$ws.Cells.Item(10, 4).Value = 598
$ws.Cells.Item(10, 5).Value = 123
$ws.Cells.Item(10, 6).Value = 4
$ws.Cells.Item(10, 7).Value = 4
$ws.Cells.Item(10, 9).Value = 614
$ws.Cells.Item(10, 10).Value = -2.605863192182412

# Row 11
$ws.Cells.Item(11, 3).Value = 251
$ws.Cells.Item(11, 4).Value = 251
$ws.Cells.Item(11, 7).Value = 6
$ws.Cells.Item(11, 9).Value = 309
$ws.Cells.Item(11, 10).Value = -18.77022653721683

# Row 12
$ws.Cells.Item(12, 2).Value = 6
$ws.Cells.Item(12, 3).Value = 435
$ws.Cells.Item(12, 4).Value = 679
$ws.Cells.Item(12, 5).Value = 92
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 8
$ws.Cells.Item(12, 8).Value = 146
$ws.Cells.Item(12, 9).Value = 789.1
$ws.Cells.Item(12, 10).Value = -13.95260423267013

# Row 13
$ws.Cells.Item(13, 3).Value = 34
$ws.Cells.Item(13, 4).Value = 36
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 9).Value = 483
$ws.Cells.Item(13, 10).Value = -92.54658385093167

# Row 14
$ws.Cells.Item(14, 2).Value = 6
$ws.Cells.Item(14, 3).Value = 624
$ws.Cells.Item(14, 4).Value = 761
$ws.Cells.Item(14, 5).Value = 61
$ws.Cells.Item(14, 6).Value = 6
$ws.Cells.Item(14, 7).Value = 8
$ws.Cells.Item(14, 8).Value = 64
$ws.Cells.Item(14, 9).Value = 719
$ws.Cells.Item(14, 10).Value = 5.841446453407517

# Row 15
$ws.Cells.Item(15, 3).Value = 201
$ws.Cells.Item(15, 4).Value = 226
$ws.Cells.Item(15, 5).Value = 24
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 1
$ws.Cells.Item(15, 9).Value = 225
$ws.Cells.Item(15, 10).Value = 0.4444444444444473

# Row 17
$ws.Cells.Item(17, 3).Value = 73
$ws.Cells.Item(17, 4).Value = 80
$ws.Cells.Item(17, 5).Value = 6
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 9).Value = 107
$ws.Cells.Item(17, 10).Value = -25.23364485981309

# Row 18
$ws.Cells.Item(18, 3).Value = 1
$ws.Cells.Item(18, 4).Value = 3
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 9).Value = 5
$ws.Cells.Item(18, 10).Value = -40

# Row 19
$ws.Cells.Item(19, 3).Value = 10
$ws.Cells.Item(19, 4).Value = 10
$ws.Cells.Item(19, 9).Value = 13
$ws.Cells.Item(19, 10).Value = -23.07692307692307

# Row 20
$ws.Cells.Item(20, 3).Value = 26
$ws.Cells.Item(20, 4).Value = 27
$ws.Cells.Item(20, 5).Value = 1
$ws.Cells.Item(20, 9).Value = 61
$ws.Cells.Item(20, 10).Value = -55.73770491803278
